# Update Name of Algo - apply updated imputed values in columns A and D
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -20.81369999999997
$ws.Range("D4").Value = -7.491000000000001
$ws.Range("A7").Value = -20.37989999999996
$ws.Range("D12").Value = -5.995799999999997
$ws.Range("A16").Value = -22.03900000000001
$ws.Range("D18").Value = -8.876099999999997
$ws.Range("D19").Value = -8.59289999999999
$ws.Range("D20").Value = -8.693599999999989
$ws.Range("A28").Value = -21.83369999999998
$ws.Range("A29").Value = -21.38279999999997
$ws.Range("D31").Value = -7.723599999999998
$ws.Range("A32").Value = -21.2247
$ws.Range("A40").Value = -20.2266
$ws.Range("D40").Value = -7.413099999999995
$ws.Range("D42").Value = -8.655399999999997
$ws.Range("D47").Value = -7.5615
$ws.Range("D48").Value = -7.294199999999997
$ws.Range("A52").Value = -22.11429999999999
$ws.Range("A57").Value = -22.53790000000002
$ws.Range("D63").Value = -6.680299999999994
$ws.Range("D64").Value = -7.066199999999993
$ws.Range("A66").Value = -21.50469999999999
$ws.Range("D76").Value = -7.553399999999995
$ws.Range("D81").Value = -8.0372
$ws.Range("D89").Value = -8.420800000000002
$ws.Range("D94").Value = -6.027899999999996
$ws.Range("A100").Value = -22.01250000000002
